$d = $word.ActiveDocument

# Locate the target paragraph: the final paragraph of the document, which
# currently holds the "This is to make it easier to understand ..." /
# ", stegosaur cannot eat ..." / " " / "Technically, location could store ..."
# runs. We replace this paragraph's whole range with three paragraphs:
#   1) a duplicate of the original paragraph (same 4 runs, same text)
#   2) a new list item: "Each item has its own location stored as an attribute"
#   3) the original paragraph, but with its runs collapsed into a single run
#      carrying new wording about travelling to items such as a corpse.

$targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$targetRange = $targetPara.Range

if ($targetRange.Text -notmatch "This is to make it easier to understand where the fruit is located") {
    throw "Target paragraph text did not match the expected content; aborting to avoid corrupting the document."
}

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPrNormal = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$rPrEastAsia = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman" w:hint="eastAsia"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$paraA = '<w:p ' + $ns + '>' +
           '<w:pPr>' + $rPrNormal + '</w:pPr>' +
           '<w:r>' + $rPrNormal + '<w:t>This is to make it easier to understand where the fruit is located as depending on where the fruit is located, different actors will interact with the fruit differently. For example</w:t></w:r>' +
           '<w:r>' + $rPrNormal + '<w:t>, stegosaur cannot eat fruit from the tree but can eat it dropped fruit.</w:t></w:r>' +
           '<w:r>' + $rPrNormal + '<w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:r>' + $rPrNormal + '<w:t>Technically, location could store dropped fruit however, it is easier to understand the code if Tree has its own dropped fruit attribute. In order to check if a location has dropped fruit under a tree and dropped fruit was stored under Location, it would involve checking the location’s ArrayList of Items as well as checking if there is a tree at that location which is more clunky than if we were to have a dropped Fruit attribute in Tree.</w:t></w:r>' +
         '</w:p>'

$paraB = '<w:p ' + $ns + '>' +
           '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + $rPrNormal + '</w:pPr>' +
           '<w:r>' + $rPrNormal + '<w:t>Each item has its own location stored as an attribute</w:t></w:r>' +
         '</w:p>'

$paraC = '<w:p ' + $ns + '>' +
           '<w:pPr>' + $rPrEastAsia + '</w:pPr>' +
           '<w:r>' + $rPrNormal + '<w:t>This is to make it easier to travel to an item such as to a corpse. Items do not travel as much as actors so the location can just be stored unlike for Actors.</w:t></w:r>' +
         '</w:p>'

$targetRange.InsertXML($paraA + $paraB + $paraC)
